$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new numeric row A5
$ws.Range("A5").Value = 127867

# Update the active selection to reflect D10 (as recorded in the saved file)
$ws.Range("D10").Select()
